$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 620 (pushes previous rows 620-699 down to 621-700)
$ws.Rows.Item(620).Insert()

# Populate the newly inserted row 620 with the new data record
$ws.Range("A620").Value = 5
$ws.Range("B620").Value = "Macroferia Regional de Talca"
$ws.Range("C620").Value = "Maule"
$ws.Range("D620").Value = 45124
$ws.Range("E620").Value = 7
$ws.Range("F620").Value = 100112043
$ws.Range("G620").Value = "Pepino ensalada"
$ws.Range("H620").Value = "Sin especificar"
$ws.Range("I620").Value = "Primera"
$ws.Range("J620").Value = 400
$ws.Range("K620").Value = 14000
$ws.Range("L620").Value = 14000
$ws.Range("M620").Value = 14000
$ws.Range("N620").Value = "$/caja 60 unidades"
$ws.Range("O620").Value = "Región de Arica y Parinacota"
$ws.Range("P620").Value = 233
$ws.Range("Q620").Value = 60
$ws.Range("R620").Value = "Hortaliza"
